{"js": "const pairs = [\n  [\"2024-09-12 Thursday\", \"2024-09-13 Friday\"],\n  [\"447\u00f78=55, 7\", \"815\u00f74=203, 3\"],\n  [\"555\u00f79=61, 6\", \"459\u00f78=57, 3\"],\n  [\"798\u00f77=114, 0\", \"345\u00f76=57, 3\"],\n  [\"583\u00f73=194, 1\", \"494\u00f78=61, 6\"],\n  [\"375\u00f73=125, 0\", \"701\u00f76=116, 5\"],\n  [\"222\u00f79=24, 6\", \"748\u00f72=374, 0\"],\n  [\"280\u00f73=93, 1\", \"613\u00f73=204, 1\"],\n  [\"944\u00f72=472, 0\", \"264\u00f77=37, 5\"],\n  [\"449\u00f74=112, 1\", \"837\u00f79=93, 0\"],\n  [\"732\u00f77=104, 4\", \"915\u00f72=457, 1\"],\n  [\"639\u00f76=106, 3\", \"888\u00f76=148, 0\"],\n  [\"588\u00f72=294, 0\", \"446\u00f74=111, 2\"],\n  [\"739\u00f74=184, 3\", \"579\u00f73=193, 0\"],\n  [\"741\u00f74=185, 1\", \"695\u00f75=139, 0\"],\n  [\"925\u00f76=154, 1\", \"787\u00f77=112, 3\"],\n  [\"724\u00f77=103, 3\", \"414\u00f72=207, 0\"],\n  [\"853\u00f77=121, 6\", \"881\u00f77=125, 6\"],\n  [\"574\u00f72=287, 0\", \"348\u00f76=58, 0\"],\n  [\"194\u00f78=24, 2\", \"264\u00f73=88, 0\"],\n  [\"512\u00f72=256, 0\", \"961\u00f77=137, 2\"],\n  [\"104\u00f76=17, 2\", \"343\u00f77=49, 0\"],\n  [\"891\u00f74=222, 3\", \"658\u00f75=131, 3\"],\n  [\"487\u00f78=60, 7\", \"108\u00f72=54, 0\"],\n  [\"500\u00f79=55, 5\", \"251\u00f72=125, 1\"],\n  [\"976\u00f78=122, 0\", \"747\u00f77=106, 5\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Not found: \" + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n# wdReplaceOne = 1  (replace a single occurrence found by the search)\n$pairs = @(\n    @(\"2024-09-12 Thursday\", \"2024-09-13 Friday\"),\n    @(\"447\u00f78=55, 7\", \"815\u00f74=203, 3\"),\n    @(\"555\u00f79=61, 6\", \"459\u00f78=57, 3\"),\n    @(\"798\u00f77=114, 0\", \"345\u00f76=57, 3\"),\n    @(\"583\u00f73=194, 1\", \"494\u00f78=61, 6\"),\n    @(\"375\u00f73=125, 0\", \"701\u00f76=116, 5\"),\n    @(\"222\u00f79=24, 6\", \"748\u00f72=374, 0\"),\n    @(\"280\u00f73=93, 1\", \"613\u00f73=204, 1\"),\n    @(\"944\u00f72=472, 0\", \"264\u00f77=37, 5\"),\n    @(\"449\u00f74=112, 1\", \"837\u00f79=93, 0\"),\n    @(\"732\u00f77=104, 4\", \"915\u00f72=457, 1\"),\n    @(\"639\u00f76=106, 3\", \"888\u00f76=148, 0\"),\n    @(\"588\u00f72=294, 0\", \"446\u00f74=111, 2\"),\n    @(\"739\u00f74=184, 3\", \"579\u00f73=193, 0\"),\n    @(\"741\u00f74=185, 1\", \"695\u00f75=139, 0\"),\n    @(\"925\u00f76=154, 1\", \"787\u00f77=112, 3\"),\n    @(\"724\u00f77=103, 3\", \"414\u00f72=207, 0\"),\n    @(\"853\u00f77=121, 6\", \"881\u00f77=125, 6\"),\n    @(\"574\u00f72=287, 0\", \"348\u00f76=58, 0\"),\n    @(\"194\u00f78=24, 2\", \"264\u00f73=88, 0\"),\n    @(\"512\u00f72=256, 0\", \"961\u00f77=137, 2\"),\n    @(\"104\u00f76=17, 2\", \"343\u00f77=49, 0\"),\n    @(\"891\u00f74=222, 3\", \"658\u00f75=131, 3\"),\n    @(\"487\u00f78=60, 7\", \"108\u00f72=54, 0\"),\n    @(\"500\u00f79=55, 5\", \"251\u00f72=125, 1\"),\n    @(\"976\u00f78=122, 0\", \"747\u00f77=106, 5\"),\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 1)\n    if (-not $found) {\n        throw \"Not found: $oldText\"\n    }\n}"}
